$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.357478060742892
$ws.Range("C2").Value = 0.1530760410253293
$ws.Range("D2").Value = 0.1194354720369049
$ws.Range("E2").Value = 0.1245758532212635
$ws.Range("F2").Value = 1.752131314218971
$ws.Range("J2").Value = 0.1612990383297617
$ws.Range("L2").Value = 0.2744277187221442
$ws.Range("M2").Value = 0.3184569435170275
$ws.Range("N2").Value = 1.721821525170089
$ws.Range("O2").Value = 4.575683049327665
$ws.Range("B3").Value = 1.279850709966922
$ws.Range("C3").Value = 0.1432450312874352
$ws.Range("D3").Value = 0.1192628499365505
$ws.Range("E3").Value = 0.1252683450860923
$ws.Range("F3").Value = 1.756729400771889
$ws.Range("J3").Value = 0.1622816523651021
$ws.Range("L3").Value = 0.2715942972574297
$ws.Range("M3").Value = 0.3060752846758632
$ws.Range("N3").Value = 1.738201905344958
$ws.Range("O3").Value = 4.589442455110543
$ws.Range("B4").Value = 1.232583085629301
$ws.Range("C4").Value = 0.137161272822425
$ws.Range("D4").Value = 0.1191855449300228
$ws.Range("E4").Value = 0.1257211810430059
$ws.Range("F4").Value = 1.760423353963681
$ws.Range("J4").Value = 0.1629171745431517
$ws.Range("L4").Value = 0.2699455850557086
$ws.Range("M4").Value = 0.2985791830747075
$ws.Range("N4").Value = 1.748849901457898
$ws.Range("O4").Value = 4.600302111763853
$ws.Range("B5").Value = 1.213422029167134
$ws.Range("C5").Value = 0.134670275636779
$ws.Range("D5").Value = 0.1191612859695432
$ws.Range("E5").Value = 0.1259126794800283
$ws.Range("F5").Value = 1.762147729475615
$ws.Range("J5").Value = 0.1631842639412078
$ws.Range("L5").Value = 0.2692967146944909
$ws.Range("M5").Value = 0.2955514316987973
$ws.Range("N5").Value = 1.753337581229914
$ws.Range("O5").Value = 4.605333945373928
$ws.Range("B6").Value = 1.21024648288298
$ws.Range("C6").Value = 0.1342559373998284
$ws.Range("D6").Value = 0.1191576961644358
$ws.Range("E6").Value = 0.1259448986474507
$ws.Range("F6").Value = 1.762447294712246
$ws.Range("J6").Value = 0.1632291040967639
$ws.Range("L6").Value = 0.269190362109299
$ws.Range("M6").Value = 0.2950503120599066
$ws.Range("N6").Value = 1.75409173059937
$ws.Range("O6").Value = 4.606206107380871
$ws.Range("B7").Value = 1.232324262673529
$ws.Range("C7").Value = 0.1371277260433175
$ws.Range("D7").Value = 0.1191851883959103
$ws.Range("E7").Value = 0.1257237354430485
$ws.Range("F7").Value = 1.76044572238289
$ws.Range("J7").Value = 0.1629207437530522
$ws.Range("L7").Value = 0.2699367409234128
$ws.Range("M7").Value = 0.2985382402327303
$ws.Range("N7").Value = 1.74890982235754
$ws.Range("O7").Value = 4.60036751728677
$ws.Range("B8").Value = 1.330630759738483
$ws.Range("C8").Value = 0.1496962338598564
$ws.Range("D8").Value = 0.1193700161124092
$ws.Range("E8").Value = 0.1248088949287585
$ws.Range("F8").Value = 1.753536082980446
$ws.Range("J8").Value = 0.1616311702708146
$ws.Range("L8").Value = 0.2734319256640276
$ws.Range("M8").Value = 0.3141658296618033
$ws.Range("N8").Value = 1.727346996114836
$ws.Range("O8").Value = 4.579926865676583
$ws.Range("B9").Value = 1.526502347138148
$ws.Range("C9").Value = 0.1739620133618018
$ws.Range("D9").Value = 0.1199588316441051
$ws.Range("E9").Value = 0.1232336408117707
$ws.Range("F9").Value = 1.746891606093897
$ws.Range("J9").Value = 0.1593571148721984
$ws.Range("L9").Value = 0.2810040988719038
$ws.Range("M9").Value = 0.3456458571932899
$ws.Range("N9").Value = 1.68974283305478
$ws.Range("O9").Value = 4.558976016853478
$ws.Range("B10").Value = 1.672245089164392
$ws.Range("C10").Value = 0.1915536162411797
$ws.Range("D10").Value = 0.1205278997622869
$ws.Range("E10").Value = 0.1222088324084538
$ws.Range("F10").Value = 1.746217339311045
$ws.Range("J10").Value = 0.1578407677485405
$ws.Range("L10").Value = 0.287000508617993
$ws.Range("M10").Value = 0.3692734304709333
$ws.Range("N10").Value = 1.664962582492919
$ws.Range("O10").Value = 4.555254460797357
$ws.Range("B11").Value = 1.738936320919834
$ws.Range("C11").Value = 0.1995043562102978
$ws.Range("D11").Value = 0.1208161000317602
$ws.Range("E11").Value = 0.1217712239271407
$ws.Range("F11").Value = 1.746823906355829
$ws.Range("J11").Value = 0.1571842709710776
$ws.Range("L11").Value = 0.2898215493165566
$ws.Range("M11").Value = 0.380128778997296
$ws.Range("N11").Value = 1.654306541597798
$ws.Range("O11").Value = 4.556097645375218
$ws.Range("B12").Value = 1.764245798480147
$ws.Range("C12").Value = 0.202507546632944
$ws.Range("D12").Value = 0.1209294228335267
$ws.Range("E12").Value = 0.1216096098853372
$ws.Range("F12").Value = 1.747184865425695
$ws.Range("J12").Value = 0.156940446719628
$ws.Range("L12").Value = 0.290903117653869
$ws.Range("M12").Value = 0.384254598846816
$ws.Range("N12").Value = 1.650359984906459
$ws.Range("O12").Value = 4.556781690891029
$ws.Range("B13").Value = 1.758792528258425
$ws.Range("C13").Value = 0.2018610952590336
$ws.Range("D13").Value = 0.1209048309679446
$ws.Range("E13").Value = 0.1216442342445996
$ws.Range("F13").Value = 1.747101289411404
$ws.Range("J13").Value = 0.1569927463617002
$ws.Range("L13").Value = 0.2906695926745329
$ws.Range("M13").Value = 0.3833653606980079
$ws.Range("N13").Value = 1.651206003743475
$ws.Range("O13").Value = 4.556618146208308
$ws.Range("B14").Value = 1.741017453675965
$ws.Range("C14").Value = 0.1997515832434544
$ws.Range("D14").Value = 0.1208253394109917
$ws.Range("E14").Value = 0.121757845761012
$ws.Range("F14").Value = 1.746850972297651
$ws.Range("J14").Value = 0.1571641157324986
$ws.Range("L14").Value = 0.2899102646248082
$ws.Range("M14").Value = 0.380467910830049
$ws.Range("N14").Value = 1.653980079298023
$ws.Range("O14").Value = 4.556146611760397
$ws.Range("B15").Value = 1.730136821615076
$ws.Range("C15").Value = 0.1984584546973451
$ws.Range("D15").Value = 0.1207771930065036
$ws.Range("E15").Value = 0.1218279695460973
$ws.Range("F15").Value = 1.746714738516403
$ws.Range("J15").Value = 0.1572697061723969
$ws.Range("L15").Value = 0.2894468833427482
$ws.Range("M15").Value = 0.3786951024844001
$ws.Range("N15").Value = 1.655690825783644
$ws.Range("O15").Value = 4.555905286154967
$ws.Range("B16").Value = 1.667894428471016
$ws.Range("C16").Value = 0.1910329631695618
$ws.Range("D16").Value = 0.1205096528104406
$ws.Range("E16").Value = 0.1222380053303009
$ws.Range("F16").Value = 1.746196070959428
$ws.Range("J16").Value = 0.1578843404151602
$ws.Range("L16").Value = 0.2868180141341412
$ws.Range("M16").Value = 0.3685661398856155
$ws.Range("N16").Value = 1.665671386456538
$ws.Range("O16").Value = 4.555250395097033
$ws.Range("B17").Value = 1.629810149942898
$ws.Range("C17").Value = 0.1864643072160277
$ws.Range("D17").Value = 0.1203530160377042
$ws.Range("E17").Value = 0.1224968614991351
$ws.Range("F17").Value = 1.746111755306018
$ws.Range("J17").Value = 0.1582699187598013
$ws.Range("L17").Value = 0.2852290938582343
$ws.Range("M17").Value = 0.3623795827904956
$ws.Range("N17").Value = 1.671952056578377
$ws.Range("O17").Value = 4.55549824659704
$ws.Range("B18").Value = 1.607942074555979
$ws.Range("C18").Value = 0.1838316726948506
$ws.Range("D18").Value = 0.1202656850031261
$ws.Range("E18").Value = 0.1226484396630836
$ws.Range("F18").Value = 1.746149226654069
$ws.Range("J18").Value = 0.1584948279738896
$ws.Range("L18").Value = 0.284323971865561
$ws.Range("M18").Value = 0.3588313299505046
$ws.Range("N18").Value = 1.675622571627642
$ws.Range("O18").Value = 4.555879514251558
$ws.Range("B19").Value = 1.600544317460844
$ws.Range("C19").Value = 0.1829394766153882
$ws.Range("D19").Value = 0.1202365915500039
$ws.Range("E19").Value = 0.1227002239412522
$ws.Range("F19").Value = 1.746176680802392
$ws.Range("J19").Value = 0.1585715171110542
$ws.Range("L19").Value = 0.2840190248484049
$ws.Range("M19").Value = 0.3576316935522215
$ws.Range("N19").Value = 1.676875313843624
$ws.Range("O19").Value = 4.55604960190999
$ws.Range("B20").Value = 1.633860468011221
$ws.Range("C20").Value = 0.1869511528948919
$ws.Range("D20").Value = 0.1203694046191401
$ws.Range("E20").Value = 0.1224690273975657
$ws.Range("F20").Value = 1.746111834010946
$ws.Range("J20").Value = 0.1582285489338839
$ws.Range("L20").Value = 0.285397328909653
$ws.Range("M20").Value = 0.363037109828241
$ws.Range("N20").Value = 1.671277461745198
$ws.Range("O20").Value = 4.555447156376317
$ws.Range("B21").Value = 1.746236943879069
$ws.Range("C21").Value = 0.200371404935737
$ws.Range("D21").Value = 0.1208485745938503
$ws.Range("E21").Value = 0.1217243641539145
$ws.Range("F21").Value = 1.746920934519537
$ws.Range("J21").Value = 0.1571136508412865
$ws.Range("L21").Value = 0.2901329375326753
$ws.Range("M21").Value = 0.3813185533891996
$ws.Range("N21").Value = 1.653162860092749
$ws.Range("O21").Value = 4.556275213204714
$ws.Range("B22").Value = 1.820001126835393
$ws.Range("C22").Value = 0.2090980564125857
$ws.Range("D22").Value = 0.121186134769502
$ws.Range("E22").Value = 0.121261569406264
$ws.Range("F22").Value = 1.748214812819171
$ws.Range("J22").Value = 0.1564128383781567
$ws.Range("L22").Value = 0.293305413134604
$ws.Range("M22").Value = 0.3933546321384469
$ws.Range("N22").Value = 1.641840664693035
$ws.Range("O22").Value = 4.558942463405515
$ws.Range("B23").Value = 1.78060301780647
$ws.Range("C23").Value = 0.2044445735252225
$ws.Range("D23").Value = 0.1210037502182999
$ws.Range("E23").Value = 0.1215063898282729
$ws.Range("F23").Value = 1.747454261913376
$ws.Range("J23").Value = 0.1567843316960298
$ws.Range("L23").Value = 0.291605148877494
$ws.Range("M23").Value = 0.3869227745812935
$ws.Range("N23").Value = 1.647836255007341
$ws.Range("O23").Value = 4.557324347926027
$ws.Range("B24").Value = 1.632029234738354
$ws.Range("C24").Value = 0.1867310688079442
$ws.Range("D24").Value = 0.1203619868627541
$ws.Range("E24").Value = 0.1224816026023685
$ws.Range("F24").Value = 1.746111530706642
$ws.Range("J24").Value = 0.1582472421543279
$ws.Range("L24").Value = 0.2853212437653809
$ws.Range("M24").Value = 0.3627398154166244
$ws.Range("N24").Value = 1.671582260183563
$ws.Range("O24").Value = 4.555469510493566
$ws.Range("B25").Value = 1.47318791948237
$ws.Range("C25").Value = 0.1674387017929462
$ws.Range("D25").Value = 0.1197754662600943
$ws.Range("E25").Value = 0.1236364558885352
$ws.Range("F25").Value = 1.747950084526238
$ws.Range("J25").Value = 0.1599451260515878
$ws.Range("L25").Value = 0.2788792273106608
$ws.Range("M25").Value = 0.3370413301114041
$ws.Range("N25").Value = 1.699415296726929
$ws.Range("O25").Value = 4.562594723743189
